# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The "K" column (G) needs its values recalculated/rewritten with the new s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 0
    4  = 1
    5  = 2
    6  = 0
    7  = 4
    8  = 1
    9  = 1
    10 = 4
    11 = 5
    12 = 3
    13 = 1
    14 = 6
    15 = 4
    16 = 5
    17 = 2
    18 = 5
    19 = 3
    20 = 3
    21 = 1
    22 = 3
    23 = 0
    24 = 4
    25 = 0
    26 = 0
    27 = 3
    28 = 7
    29 = 4
    30 = 2
    31 = 1
    32 = 0
    33 = 1
    34 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
